$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new row 24 date cell uses the same custom date number format as the
# existing date column (D2:D23), matching style index "s=2" in the template.
$ws.Cells.Item(24, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 24 is a brand-new data row; populate the columns that are constant
# across every row in this sheet.
$ws.Cells.Item(24, 1).Value = 7
$ws.Cells.Item(24, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(24, 3).Value = "Ñuble"
$ws.Cells.Item(24, 5).Value = 16
$ws.Cells.Item(24, 6).Value = "Fruta"
$ws.Cells.Item(24, 7).Value = 100107
$ws.Cells.Item(24, 8).Value = "Otros"
$ws.Cells.Item(24, 9).Value = 100107002
$ws.Cells.Item(24, 10).Value = "Chirimoya"
$ws.Cells.Item(24, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(24, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(24, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(24, 20).Value = 10

# Weekly refresh of the Chirimoya price rows: dates, quality grade, volumes
# and min/max/weighted prices for rows 2-24.
$ws.Cells.Item(2, 4).Value = "9/13/2021"
$ws.Cells.Item(2, 12).Value = "Primera"
$ws.Cells.Item(2, 13).Value = 60
$ws.Cells.Item(2, 14).Value = 21000
$ws.Cells.Item(2, 15).Value = 22000
$ws.Cells.Item(2, 16).Value = 21500
$ws.Cells.Item(2, 19).Value = 2150

$ws.Cells.Item(3, 4).Value = "10/5/2022"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 120
$ws.Cells.Item(3, 14).Value = 25000
$ws.Cells.Item(3, 15).Value = 26000
$ws.Cells.Item(3, 16).Value = 25500
$ws.Cells.Item(3, 19).Value = 2550

$ws.Cells.Item(4, 4).Value = "11/3/2022"
$ws.Cells.Item(4, 12).Value = "Especial"
$ws.Cells.Item(4, 13).Value = 60
$ws.Cells.Item(4, 14).Value = 26000
$ws.Cells.Item(4, 15).Value = 26000
$ws.Cells.Item(4, 16).Value = 26000
$ws.Cells.Item(4, 19).Value = 2600

$ws.Cells.Item(5, 4).Value = "9/22/2023"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 30
$ws.Cells.Item(5, 14).Value = 21000
$ws.Cells.Item(5, 15).Value = 21000
$ws.Cells.Item(5, 16).Value = 21000
$ws.Cells.Item(5, 19).Value = 2100

$ws.Cells.Item(6, 4).Value = "10/18/2021"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 30
$ws.Cells.Item(6, 14).Value = 23000
$ws.Cells.Item(6, 15).Value = 24000
$ws.Cells.Item(6, 16).Value = 23500
$ws.Cells.Item(6, 19).Value = 2350

$ws.Cells.Item(7, 4).Value = "9/20/2023"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 50
$ws.Cells.Item(7, 14).Value = 22000
$ws.Cells.Item(7, 15).Value = 22000
$ws.Cells.Item(7, 16).Value = 22000
$ws.Cells.Item(7, 19).Value = 2200

$ws.Cells.Item(8, 4).Value = "11/9/2022"
$ws.Cells.Item(8, 12).Value = "Especial"
$ws.Cells.Item(8, 13).Value = 30
$ws.Cells.Item(8, 14).Value = 25000
$ws.Cells.Item(8, 15).Value = 25000
$ws.Cells.Item(8, 16).Value = 25000
$ws.Cells.Item(8, 19).Value = 2500

$ws.Cells.Item(9, 4).Value = "11/9/2022"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 80
$ws.Cells.Item(9, 14).Value = 23000
$ws.Cells.Item(9, 15).Value = 24000
$ws.Cells.Item(9, 16).Value = 23500
$ws.Cells.Item(9, 19).Value = 2350

$ws.Cells.Item(10, 4).Value = "9/27/2023"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 50
$ws.Cells.Item(10, 14).Value = 23000
$ws.Cells.Item(10, 15).Value = 23000
$ws.Cells.Item(10, 16).Value = 23000
$ws.Cells.Item(10, 19).Value = 2300

$ws.Cells.Item(11, 4).Value = "9/8/2021"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 60
$ws.Cells.Item(11, 14).Value = 21000
$ws.Cells.Item(11, 15).Value = 22000
$ws.Cells.Item(11, 16).Value = 21500
$ws.Cells.Item(11, 19).Value = 2150

$ws.Cells.Item(12, 4).Value = "9/4/2023"
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 50
$ws.Cells.Item(12, 14).Value = 22000
$ws.Cells.Item(12, 15).Value = 22000
$ws.Cells.Item(12, 16).Value = 22000
$ws.Cells.Item(12, 19).Value = 2200

$ws.Cells.Item(13, 4).Value = "9/21/2021"
$ws.Cells.Item(13, 12).Value = "Especial"
$ws.Cells.Item(13, 13).Value = 60
$ws.Cells.Item(13, 14).Value = 31000
$ws.Cells.Item(13, 15).Value = 32000
$ws.Cells.Item(13, 16).Value = 31500
$ws.Cells.Item(13, 19).Value = 3150

$ws.Cells.Item(14, 4).Value = "9/21/2021"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 30
$ws.Cells.Item(14, 14).Value = 30000
$ws.Cells.Item(14, 15).Value = 30000
$ws.Cells.Item(14, 16).Value = 30000
$ws.Cells.Item(14, 19).Value = 3000

$ws.Cells.Item(15, 4).Value = "9/9/2021"
$ws.Cells.Item(15, 12).Value = "Primera"
$ws.Cells.Item(15, 13).Value = 60
$ws.Cells.Item(15, 14).Value = 21000
$ws.Cells.Item(15, 15).Value = 22000
$ws.Cells.Item(15, 16).Value = 21500
$ws.Cells.Item(15, 19).Value = 2150

$ws.Cells.Item(16, 4).Value = "9/11/2023"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 40
$ws.Cells.Item(16, 14).Value = 22000
$ws.Cells.Item(16, 15).Value = 22000
$ws.Cells.Item(16, 16).Value = 22000
$ws.Cells.Item(16, 19).Value = 2200

$ws.Cells.Item(17, 4).Value = "9/22/2021"
$ws.Cells.Item(17, 12).Value = "Especial"
$ws.Cells.Item(17, 13).Value = 60
$ws.Cells.Item(17, 14).Value = 31000
$ws.Cells.Item(17, 15).Value = 32000
$ws.Cells.Item(17, 16).Value = 31500
$ws.Cells.Item(17, 19).Value = 3150

$ws.Cells.Item(18, 4).Value = "9/22/2021"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 30
$ws.Cells.Item(18, 14).Value = 30000
$ws.Cells.Item(18, 15).Value = 30000
$ws.Cells.Item(18, 16).Value = 30000
$ws.Cells.Item(18, 19).Value = 3000

$ws.Cells.Item(19, 4).Value = "10/14/2022"
$ws.Cells.Item(19, 12).Value = "Especial"
$ws.Cells.Item(19, 13).Value = 60
$ws.Cells.Item(19, 14).Value = 24000
$ws.Cells.Item(19, 15).Value = 25000
$ws.Cells.Item(19, 16).Value = 24500
$ws.Cells.Item(19, 19).Value = 2450

$ws.Cells.Item(20, 4).Value = "10/14/2022"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 120
$ws.Cells.Item(20, 14).Value = 21000
$ws.Cells.Item(20, 15).Value = 22000
$ws.Cells.Item(20, 16).Value = 21500
$ws.Cells.Item(20, 19).Value = 2150

$ws.Cells.Item(21, 4).Value = "10/7/2022"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 60
$ws.Cells.Item(21, 14).Value = 23000
$ws.Cells.Item(21, 15).Value = 24000
$ws.Cells.Item(21, 16).Value = 23500
$ws.Cells.Item(21, 19).Value = 2350

$ws.Cells.Item(22, 4).Value = "9/7/2021"
$ws.Cells.Item(22, 12).Value = "Primera"
$ws.Cells.Item(22, 13).Value = 60
$ws.Cells.Item(22, 14).Value = 21000
$ws.Cells.Item(22, 15).Value = 22000
$ws.Cells.Item(22, 16).Value = 21500
$ws.Cells.Item(22, 19).Value = 2150

$ws.Cells.Item(23, 4).Value = "9/25/2023"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 80
$ws.Cells.Item(23, 14).Value = 22000
$ws.Cells.Item(23, 15).Value = 22000
$ws.Cells.Item(23, 16).Value = 22000
$ws.Cells.Item(23, 19).Value = 2200

$ws.Cells.Item(24, 4).Value = "9/7/2023"
$ws.Cells.Item(24, 12).Value = "Primera"
$ws.Cells.Item(24, 13).Value = 30
$ws.Cells.Item(24, 14).Value = 22000
$ws.Cells.Item(24, 15).Value = 22000
$ws.Cells.Item(24, 16).Value = 22000
$ws.Cells.Item(24, 19).Value = 2200

